$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension implicitly handled by writing cells; extend used range to row 12

# Row 2
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = '5862'
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = '5/26/2025'
$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = 'CHARCAS 3715'
$ws.Cells.Item(2, 4).Value = 14
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '806976061'
$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = 'Optical Power'
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = 'Pendiente'
$ws.Cells.Item(2, 8).NumberFormat = "@"
$ws.Cells.Item(2, 8).Value = 'Cable en panza Cable cortado'
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = '{"direccionesNormalizadas": [{"altura": 3715, "cod_calle": 3219, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.417181", "y": "-34.588033"}, "direccion": "CHARCAS 3715, CABA", "nombre_calle": "CHARCAS", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(2, 11).Value = -58.417181
$ws.Cells.Item(2, 12).Value = -34.588033
$ws.Cells.Item(2, 13).NumberFormat = "@"
$ws.Cells.Item(2, 13).Value = 'Palermo'
$ws.Cells.Item(2, 14).NumberFormat = "@"
$ws.Cells.Item(2, 14).Value = 'Capital Sur'

# Row 3
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = '5894'
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = '5/27/2025'
$ws.Cells.Item(3, 3).NumberFormat = "@"
$ws.Cells.Item(3, 3).Value = 'ALBARELLOS AV. 3100'
$ws.Cells.Item(3, 4).Value = 12
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '807045500'
$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 6).Value = 'Optical Power'
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = 'Pendiente'
$ws.Cells.Item(3, 8).NumberFormat = "@"
$ws.Cells.Item(3, 8).Value = 'No coinciden las fotos cargadas en el form con el reclamo original'
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).NumberFormat = "@"
$ws.Cells.Item(3, 10).Value = '{"direccionesNormalizadas": [{"altura": 3100, "cod_calle": 1029, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.512533", "y": "-34.579243"}, "direccion": "ALBARELLOS AV. 3100, CABA", "nombre_calle": "ALBARELLOS AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(3, 11).Value = -58.512533
$ws.Cells.Item(3, 12).Value = -34.579243
$ws.Cells.Item(3, 13).NumberFormat = "@"
$ws.Cells.Item(3, 13).Value = 'Paternal'
$ws.Cells.Item(3, 14).NumberFormat = "@"
$ws.Cells.Item(3, 14).Value = 'Capital Norte'

# Row 4
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = '4238'
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = '6/2/2025'
$ws.Cells.Item(4, 3).NumberFormat = "@"
$ws.Cells.Item(4, 3).Value = 'GUATEMALA 5527'
$ws.Cells.Item(4, 4).Value = 14
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '807130137'
$ws.Cells.Item(4, 6).NumberFormat = "@"
$ws.Cells.Item(4, 6).Value = 'Optical Power'
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = 'Pendiente'
$ws.Cells.Item(4, 8).NumberFormat = "@"
$ws.Cells.Item(4, 8).Value = 'Las fotos que cargaron en el form no coinciden con las originales'
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).NumberFormat = "@"
$ws.Cells.Item(4, 10).Value = '{"direccionesNormalizadas": [{"altura": 5527, "cod_calle": 7093, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.431302", "y": "-34.580805"}, "direccion": "GUATEMALA 5527, CABA", "nombre_calle": "GUATEMALA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(4, 11).Value = -58.431302
$ws.Cells.Item(4, 12).Value = -34.580805
$ws.Cells.Item(4, 13).NumberFormat = "@"
$ws.Cells.Item(4, 13).Value = 'Palermo'
$ws.Cells.Item(4, 14).NumberFormat = "@"
$ws.Cells.Item(4, 14).Value = 'Capital Sur'

# Row 5
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = '6262'
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = '6/25/2025'
$ws.Cells.Item(5, 3).NumberFormat = "@"
$ws.Cells.Item(5, 3).Value = 'MIGUELETES 1330'
$ws.Cells.Item(5, 4).Value = 14
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '807789707'
$ws.Cells.Item(5, 6).NumberFormat = "@"
$ws.Cells.Item(5, 6).Value = 'Optical Power'
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = 'Pendiente'
$ws.Cells.Item(5, 8).NumberFormat = "@"
$ws.Cells.Item(5, 8).Value = 'Cables en panza'
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).NumberFormat = "@"
$ws.Cells.Item(5, 10).Value = '{"direccionesNormalizadas": [{"altura": 1330, "cod_calle": 13079, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.440291", "y": "-34.562841"}, "direccion": "MIGUELETES 1330, CABA", "nombre_calle": "MIGUELETES", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(5, 11).Value = -58.440291
$ws.Cells.Item(5, 12).Value = -34.562841
$ws.Cells.Item(5, 13).NumberFormat = "@"
$ws.Cells.Item(5, 13).Value = 'Colegiales'
$ws.Cells.Item(5, 14).NumberFormat = "@"
$ws.Cells.Item(5, 14).Value = 'Capital Norte'

# Row 6
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = '6417'
$ws.Cells.Item(6, 2).NumberFormat = "@"
$ws.Cells.Item(6, 2).Value = '7/15/2025'
$ws.Cells.Item(6, 3).NumberFormat = "@"
$ws.Cells.Item(6, 3).Value = 'NUMANCIA 436'
$ws.Cells.Item(6, 4).Value = 6
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '808373678'
$ws.Cells.Item(6, 6).NumberFormat = "@"
$ws.Cells.Item(6, 6).Value = 'Optical Power'
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = 'Pendiente'
$ws.Cells.Item(6, 8).NumberFormat = "@"
$ws.Cells.Item(6, 8).Value = 'Tendido a baja altura'
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(6, 10).NumberFormat = "@"
$ws.Cells.Item(6, 10).Value = '{"direccionesNormalizadas": [{"altura": 436, "cod_calle": 14030, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.432659", "y": "-34.611358"}, "direccion": "NUMANCIA 436, CABA", "nombre_calle": "NUMANCIA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(6, 11).Value = -58.432659
$ws.Cells.Item(6, 12).Value = -34.611358
$ws.Cells.Item(6, 13).NumberFormat = "@"
$ws.Cells.Item(6, 13).Value = 'Almagro'
$ws.Cells.Item(6, 14).NumberFormat = "@"
$ws.Cells.Item(6, 14).Value = 'Capital Sur'

# Row 7
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = '6520'
$ws.Cells.Item(7, 2).NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = '7/28/2025'
$ws.Cells.Item(7, 3).NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = 'CAMPANA 380'
$ws.Cells.Item(7, 4).Value = 10
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = 'ICD30242530 '
$ws.Cells.Item(7, 6).NumberFormat = "@"
$ws.Cells.Item(7, 6).Value = 'Optical Power'
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = 'Pendiente'
$ws.Cells.Item(7, 8).NumberFormat = "@"
$ws.Cells.Item(7, 8).Value = 'Tendido a baja altura'
$ws.Cells.Item(7, 9).Value = 1
$ws.Cells.Item(7, 10).NumberFormat = "@"
$ws.Cells.Item(7, 10).Value = '{"direccionesNormalizadas": [{"altura": 380, "cod_calle": 3039, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.476471", "y": "-34.628097"}, "direccion": "CAMPANA 380, CABA", "nombre_calle": "CAMPANA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(7, 11).Value = -58.476471
$ws.Cells.Item(7, 12).Value = -34.628097
$ws.Cells.Item(7, 13).NumberFormat = "@"
$ws.Cells.Item(7, 13).Value = 'Devoto'
$ws.Cells.Item(7, 14).NumberFormat = "@"
$ws.Cells.Item(7, 14).Value = 'Capital Norte'

# Row 8
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = '6537'
$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = '7/29/2025'
$ws.Cells.Item(8, 3).NumberFormat = "@"
$ws.Cells.Item(8, 3).Value = 'CHIVILCOY 452'
$ws.Cells.Item(8, 4).Value = 10
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = 'ICD30242453 '
$ws.Cells.Item(8, 6).NumberFormat = "@"
$ws.Cells.Item(8, 6).Value = 'Optical Power'
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = 'Pendiente'
$ws.Cells.Item(8, 8).NumberFormat = "@"
$ws.Cells.Item(8, 8).Value = 'Caja de empalme colgando'
$ws.Cells.Item(8, 9).Value = 1
$ws.Cells.Item(8, 10).NumberFormat = "@"
$ws.Cells.Item(8, 10).Value = '{"direccionesNormalizadas": [{"altura": 452, "cod_calle": 3258, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.484008", "y": "-34.628912"}, "direccion": "CHIVILCOY 452, CABA", "nombre_calle": "CHIVILCOY", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(8, 11).Value = -58.484008
$ws.Cells.Item(8, 12).Value = -34.628912
$ws.Cells.Item(8, 13).NumberFormat = "@"
$ws.Cells.Item(8, 13).Value = 'Devoto'
$ws.Cells.Item(8, 14).NumberFormat = "@"
$ws.Cells.Item(8, 14).Value = 'Capital Norte'

# Row 9
$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = '6540'
$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = '7/29/2025'
$ws.Cells.Item(9, 3).NumberFormat = "@"
$ws.Cells.Item(9, 3).Value = 'CUENCA 311'
$ws.Cells.Item(9, 4).Value = 7
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = 'ICD30242102 '
$ws.Cells.Item(9, 6).NumberFormat = "@"
$ws.Cells.Item(9, 6).Value = 'Optical Power'
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = 'Pendiente'
$ws.Cells.Item(9, 8).NumberFormat = "@"
$ws.Cells.Item(9, 8).Value = 'Tendido a baja altura'
$ws.Cells.Item(9, 9).Value = 1
$ws.Cells.Item(9, 10).NumberFormat = "@"
$ws.Cells.Item(9, 10).Value = '{"direccionesNormalizadas": [{"altura": 311, "cod_calle": 3200, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.475030", "y": "-34.628307"}, "direccion": "CUENCA 311, CABA", "nombre_calle": "CUENCA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(9, 11).Value = -58.47503
$ws.Cells.Item(9, 12).Value = -34.628307
$ws.Cells.Item(9, 13).NumberFormat = "@"
$ws.Cells.Item(9, 13).Value = 'Devoto'
$ws.Cells.Item(9, 14).NumberFormat = "@"
$ws.Cells.Item(9, 14).Value = 'Capital Norte'

# Row 10
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = '6557'
$ws.Cells.Item(10, 2).NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = '8/4/2025'
$ws.Cells.Item(10, 3).NumberFormat = "@"
$ws.Cells.Item(10, 3).Value = 'ALBERDI, JUAN BAUTISTA AV. 1091'
$ws.Cells.Item(10, 4).Value = 6
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = ''
$ws.Cells.Item(10, 6).NumberFormat = "@"
$ws.Cells.Item(10, 6).Value = 'Optical Power'
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = 'Pendiente'
$ws.Cells.Item(10, 8).NumberFormat = "@"
$ws.Cells.Item(10, 8).Value = 'Cables a baja altura'
$ws.Cells.Item(10, 9).Value = 1
$ws.Cells.Item(10, 10).NumberFormat = "@"
$ws.Cells.Item(10, 10).Value = '{"direccionesNormalizadas": [{"altura": 1091, "cod_calle": 1033, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.443538", "y": "-34.622890"}, "direccion": "ALBERDI, JUAN BAUTISTA AV. 1091, CABA", "nombre_calle": "ALBERDI, JUAN BAUTISTA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(10, 11).Value = -58.443538
$ws.Cells.Item(10, 12).Value = -34.62289
$ws.Cells.Item(10, 13).NumberFormat = "@"
$ws.Cells.Item(10, 13).Value = 'Boedo'
$ws.Cells.Item(10, 14).NumberFormat = "@"
$ws.Cells.Item(10, 14).Value = 'Capital Sur'

# Row 11
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = '6193'
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = '8/4/2025'
$ws.Cells.Item(11, 3).NumberFormat = "@"
$ws.Cells.Item(11, 3).Value = 'POLA 591'
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = ''
$ws.Cells.Item(11, 6).NumberFormat = "@"
$ws.Cells.Item(11, 6).Value = 'Optical Power'
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = 'Pendiente'
$ws.Cells.Item(11, 8).NumberFormat = "@"
$ws.Cells.Item(11, 8).Value = 'Cable colgando y enrollado en arbol'
$ws.Cells.Item(11, 9).Value = 1
$ws.Cells.Item(11, 10).NumberFormat = "@"
$ws.Cells.Item(11, 10).Value = '{"direccionesNormalizadas": [{"altura": 591, "cod_calle": 17105, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.507385", "y": "-34.644479"}, "direccion": "POLA 591, CABA", "nombre_calle": "POLA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(11, 11).Value = -58.507385
$ws.Cells.Item(11, 12).Value = -34.644479
$ws.Cells.Item(11, 13).NumberFormat = "@"
$ws.Cells.Item(11, 13).Value = 'Devoto'
$ws.Cells.Item(11, 14).NumberFormat = "@"
$ws.Cells.Item(11, 14).Value = 'Capital Norte'

# Row 12
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = '6273'
$ws.Cells.Item(12, 2).NumberFormat = "@"
$ws.Cells.Item(12, 2).Value = '8/4/2025'
$ws.Cells.Item(12, 3).NumberFormat = "@"
$ws.Cells.Item(12, 3).Value = 'ARGERICH 516'
$ws.Cells.Item(12, 4).Value = 7
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = ''
$ws.Cells.Item(12, 6).NumberFormat = "@"
$ws.Cells.Item(12, 6).Value = 'Optical Power'
$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = 'Pendiente'
$ws.Cells.Item(12, 8).NumberFormat = "@"
$ws.Cells.Item(12, 8).Value = 'Tendido a baja altura'
$ws.Cells.Item(12, 9).Value = 1
$ws.Cells.Item(12, 10).NumberFormat = "@"
$ws.Cells.Item(12, 10).Value = '{"direccionesNormalizadas": [{"altura": 516, "cod_calle": 1110, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.473288", "y": "-34.626689"}, "direccion": "ARGERICH 516, CABA", "nombre_calle": "ARGERICH", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(12, 11).Value = -58.473288
$ws.Cells.Item(12, 12).Value = -34.626689
$ws.Cells.Item(12, 13).NumberFormat = "@"
$ws.Cells.Item(12, 13).Value = 'Devoto'
$ws.Cells.Item(12, 14).NumberFormat = "@"
$ws.Cells.Item(12, 14).Value = 'Capital Norte'
